$p = $ppt.ActivePresentation

# The deck currently has:
#   Slide 2 = "Code Created - Technical"
#   Slide 3 = "Code Created - Solution" (includes the architecture picture)
# They need to swap places so the Solution slide comes before the Technical
# slide (Slide 2 = Solution, Slide 3 = Technical).
$s3 = $p.Slides.Item(3)
$s3.MoveTo(2)
